$d = $word.ActiveDocument

$d.Content.Find.Execute("170÷5=34, 0", $true, $false, $false, $false, $false, $true, 1, $false, "643÷5=128, 3", 2)
$d.Content.Find.Execute("160÷3=53, 1", $true, $false, $false, $false, $false, $true, 1, $false, "233÷6=38, 5", 2)
$d.Content.Find.Execute("934÷9=103, 7", $true, $false, $false, $false, $false, $true, 1, $false, "283÷2=141, 1", 2)
$d.Content.Find.Execute("963÷8=120, 3", $true, $false, $false, $false, $false, $true, 1, $false, "765÷4=191, 1", 2)
$d.Content.Find.Execute("615÷5=123, 0", $true, $false, $false, $false, $false, $true, 1, $false, "436÷8=54, 4", 2)
$d.Content.Find.Execute("156÷5=31, 1", $true, $false, $false, $false, $false, $true, 1, $false, "403÷3=134, 1", 2)
$d.Content.Find.Execute("430÷8=53, 6", $true, $false, $false, $false, $false, $true, 1, $false, "462÷3=154, 0", 2)
$d.Content.Find.Execute("132÷4=33, 0", $true, $false, $false, $false, $false, $true, 1, $false, "395÷9=43, 8", 2)
$d.Content.Find.Execute("478÷9=53, 1", $true, $false, $false, $false, $false, $true, 1, $false, "354÷6=59, 0", 2)
$d.Content.Find.Execute("689÷9=76, 5", $true, $false, $false, $false, $false, $true, 1, $false, "624÷3=208, 0", 2)
$d.Content.Find.Execute("481÷6=80, 1", $true, $false, $false, $false, $false, $true, 1, $false, "824÷5=164, 4", 2)
$d.Content.Find.Execute("170÷9=18, 8", $true, $false, $false, $false, $false, $true, 1, $false, "232÷2=116, 0", 2)
$d.Content.Find.Execute("752÷9=83, 5", $true, $false, $false, $false, $false, $true, 1, $false, "338÷8=42, 2", 2)
$d.Content.Find.Execute("280÷4=70, 0", $true, $false, $false, $false, $false, $true, 1, $false, "894÷5=178, 4", 2)
$d.Content.Find.Execute("707÷3=235, 2", $true, $false, $false, $false, $false, $true, 1, $false, "184÷7=26, 2", 2)
$d.Content.Find.Execute("632÷4=158, 0", $true, $false, $false, $false, $false, $true, 1, $false, "403÷6=67, 1", 2)
$d.Content.Find.Execute("350÷7=50, 0", $true, $false, $false, $false, $false, $true, 1, $false, "633÷5=126, 3", 2)
$d.Content.Find.Execute("585÷8=73, 1", $true, $false, $false, $false, $false, $true, 1, $false, "615÷8=76, 7", 2)
$d.Content.Find.Execute("399÷3=133, 0", $true, $false, $false, $false, $false, $true, 1, $false, "709÷8=88, 5", 2)
$d.Content.Find.Execute("381÷5=76, 1", $true, $false, $false, $false, $false, $true, 1, $false, "250÷2=125, 0", 2)
$d.Content.Find.Execute("943÷9=104, 7", $true, $false, $false, $false, $false, $true, 1, $false, "932÷9=103, 5", 2)
$d.Content.Find.Execute("228÷8=28, 4", $true, $false, $false, $false, $false, $true, 1, $false, "218÷6=36, 2", 2)
$d.Content.Find.Execute("333÷8=41, 5", $true, $false, $false, $false, $false, $true, 1, $false, "881÷2=440, 1", 2)
$d.Content.Find.Execute("230÷9=25, 5", $true, $false, $false, $false, $false, $true, 1, $false, "819÷2=409, 1", 2)
$d.Content.Find.Execute("949÷8=118, 5", $true, $false, $false, $false, $false, $true, 1, $false, "468÷4=117, 0", 2)
